$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fifa_world_cup_2018_matches")

# Round-of-16 matches (rows 50-57): mark as completed and fill in results
$ws.Range("D50").Value = "completed"
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 3

$ws.Range("D51").Value = "completed"
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 1

$ws.Range("D52").Value = "completed"
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 1

$ws.Range("D53").Value = "completed"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 1

$ws.Range("D54").Value = "completed"
$ws.Range("G54").Value = 2

$ws.Range("D55").Value = "completed"
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 2

$ws.Range("D56").Value = "completed"
$ws.Range("G56").Value = 1

$ws.Range("D57").Value = "completed"
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 1

# Quarter-final matchups (rows 58-61): fill in the teams that advanced
$ws.Range("E58").Value = "France"
$ws.Range("F58").Value = "Uruguay"

$ws.Range("E59").Value = "Belgium"
$ws.Range("F59").Value = "Brazil"

$ws.Range("E60").Value = "England"
$ws.Range("F60").Value = "Sweden"

$ws.Range("E61").Value = "Russia"
$ws.Range("F61").Value = "Croatia"

# Update the active selection to match where the edit left off
$ws.Range("F59").Select() | Out-Null
